# Update coStatements text for the course rows (data cleanup: commas/dashes -> pipe/hash delimiters)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = "good theory#2| better lab#3"
$ws.Range("I3").Value = "learn packet tracer#2| implement TCP#1"
$ws.Range("I4").Value = "good theory#2| better lab#3"

# Move the active selection to L5, matching the editor's last cursor position
$null = $ws.Range("L5").Select()
